$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record as row 215, pushing the existing
# rows 215:250 down to 216:251 (dimension grows from A1:R250 to A1:R251).
$ws.Rows("215:215").Insert()

$ws.Range("A215").Value = 6
$ws.Range("B215").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C215").Value = "Metropolitana"
$ws.Range("D215").Value = 44476
$ws.Range("E215").Value = 13
$ws.Range("F215").Value = 100112052
$ws.Range("G215").Value = "Albahaca"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 40
$ws.Range("K215").Value = 4000
$ws.Range("L215").Value = 5000
$ws.Range("M215").Value = 4575
$ws.Range("N215").Value = "$/paquete"
$ws.Range("O215").Value = "Región de Arica y Parinacota"
$ws.Range("P215").Value = 4575
$ws.Range("Q215").Value = 1
$ws.Range("R215").Value = "Hortaliza"
